$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 45, shifting existing rows 45-53 down to 46-54.
$ws.Rows.Item(45).Insert()

# Populate the newly inserted row 45 with the new data record.
$ws.Range("A45").Value = 5
$ws.Range("B45").Value = "Macroferia Regional de Talca"
$ws.Range("C45").Value = "Maule"
$ws.Range("D45").Value = 44522
$ws.Range("E45").Value = 7
$ws.Range("F45").Value = 300000000
$ws.Range("G45").Value = "Espárragos"
$ws.Range("H45").Value = "Verde"
$ws.Range("I45").Value = "Primera"
$ws.Range("J45").Value = 3000
$ws.Range("K45").Value = 1200
$ws.Range("L45").Value = 1200
$ws.Range("M45").Value = 1200
$ws.Range("N45").Value = "`$/kilo"
$ws.Range("O45").Value = "Provincia de Linares"
$ws.Range("P45").Value = 1200
$ws.Range("Q45").Value = 1
$ws.Range("R45").Value = "Hortaliza"
